$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Icons")

# Update the "New message" icon URL (row 4, column B) to the new note icon.
$ws.Range("B4").Value = "http://www.flaticon.com/free-icon/note-black-paper-with-text-lines_34074"

# Insert three new rows above the "Doorbell 1" row (old row 28) to make room for
# two new wallpaper-selection categories; the third inserted row is left blank.
$ws.Rows("28:30").Insert()

# Fill in the two new rows: Landscapes and Holidays.
$ws.Range("A28").Value = "Landscapes"
$ws.Range("B28").Value = "http://www.flaticon.com/free-icon/landscape_92744"
$ws.Range("A29").Value = "Holidays"
$ws.Range("B29").Value = "http://www.flaticon.com/free-icon/fireworks_108980"

# Update the sheet view: scroll down a bit and leave the selection on B29.
$ws.Activate()
$ws.Range("B29").Select()
